$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New long prompt text (with updated adjacency matrix) for cell A2
$newPrompt = @'
 Given is the adjacency matrix for a unweighted undirected graph containing 16 nodes labelled A to P. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   

Consider some examples

Example 1: what is the shortest path from node A to node P?
   A B C D E F G H I J K L M N O P
 A 0 1 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 1 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 1 0 0 1 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 0 0 0 1 0 0 0 0 0 0 0 0
 E 1 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0
 F 0 1 0 0 1 0 1 0 0 1 0 0 0 0 0 0
 G 0 0 1 0 0 1 0 1 0 0 1 0 0 0 0 0
 H 0 0 0 1 0 0 1 0 0 0 0 1 0 0 0 0
 I 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0
 J 0 0 0 0 0 1 0 0 0 0 1 0 0 0 0 0
 K 0 0 0 0 0 0 1 0 0 1 0 1 0 0 1 0
 L 0 0 0 0 0 0 0 1 0 0 1 0 0 0 0 1
 M 0 0 0 0 0 0 0 0 1 0 0 0 0 1 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0
 O 0 0 0 0 0 0 0 0 0 0 1 0 0 1 0 1
 P 0 0 0 0 0 0 0 0 0 0 0 1 0 0 1 0

Solution: A -> E -> F -> G -> H -> L -> P
        

Example 2: what is the shortest path from node A to node I?
   A B C D E F G H I
 A 0 1 0 1 0 0 0 0 0
 B 1 0 1 0 1 0 0 0 0
 C 0 1 0 0 0 1 0 0 0
 D 1 0 0 0 1 0 1 0 0
 E 0 1 0 1 0 1 0 1 0
 F 0 0 1 0 1 0 0 0 1
 G 0 0 0 1 0 0 0 1 0
 H 0 0 0 0 1 0 1 0 1
 I 0 0 0 0 0 1 0 1 0

Solution: A -> D -> E -> F -> I
        

Example 3: what is the shortest path from node A to node I?
   A B C D E F G H I
 A 0 1 0 1 0 0 0 0 0
 B 1 0 1 0 0 0 0 0 0
 C 0 1 0 0 0 1 0 0 0
 D 1 0 0 0 1 0 1 0 0
 E 0 0 0 1 0 1 0 1 0
 F 0 0 1 0 1 0 0 0 1
 G 0 0 0 1 0 0 0 1 0
 H 0 0 0 0 1 0 1 0 1
 I 0 0 0 0 0 1 0 1 0

Solution: A -> D -> E -> F -> I
        
 Given these examples, answer the following quesiton.

what is the shortest path from node A to node P?

   A B C D E F G H I J K L M N O P
 A 0 1 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 1 0 0 1 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 0 0 0 1 0 0 0 0 0 0 0 0
 E 1 0 0 0 0 1 0 0 1 0 0 0 0 0 0 0
 F 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0
 G 0 0 1 0 0 1 0 1 0 0 1 0 0 0 0 0
 H 0 0 0 1 0 0 1 0 0 0 0 1 0 0 0 0
 I 0 0 0 0 1 0 0 0 0 1 0 0 1 0 0 0
 J 0 0 0 0 0 0 0 0 1 0 0 0 0 1 0 0
 K 0 0 0 0 0 0 1 0 0 0 0 1 0 0 0 0
 L 0 0 0 0 0 0 0 1 0 0 1 0 0 0 0 1
 M 0 0 0 0 0 0 0 0 1 0 0 0 0 1 0 0
 N 0 0 0 0 0 0 0 0 0 1 0 0 1 0 1 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1
 P 0 0 0 0 0 0 0 0 0 0 0 1 0 0 1 0
    
'@

# --- Update sheet1 (o_10): add column E header + row2 values, update existing text ---
$ws1.Range("A2").Value = $newPrompt
$ws1.Range("B2").Value = "A -> E -> F -> G -> H -> L -> P"
$ws1.Range("C2").Value = "The shortest path from node A to node P is A -> E -> F -> G -> H -> L -> P."
$ws1.Range("D2").Value = "Correct"

# Setting the long multi-line prompt text can trigger Excel's row-height
# autofit; restore the default row height so the row stays unchanged.
$ws1.Rows.Item(2).RowHeight = 15

# Add header E1, formatted like D1 (bold, centered, bordered)
$ws1.Range("D1").Copy($ws1.Range("E1"))
$ws1.Range("E1").Value = "evaluator_partial_correctness"
$ws1.Range("E2").Value = "Output: 7/7"

# --- Add new worksheets o_20 and o_20_jumbled, positioned after o_10 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "o_20"
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "o_20_jumbled"

# Copy header row (with formatting) from sheet1 into the new sheets
$ws1.Range("A1:E1").Copy($ws2.Range("A1:E1"))
$ws1.Range("A1:E1").Copy($ws3.Range("A1:E1"))

$ws1.Activate()
